$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Email cell (B2) to a valid address for the combined script
$ws.Range("B2").Value = "yu.li9@hpe.com"

# Move the active selection (matches the saved cursor position in the file)
$ws.Range("B8").Select() | Out-Null
